# progression.xlsx - "newer version of spreadsheet"
#
# This script reproduces (on the Sheet1 / Representations workbook) the
# content edits described by the commit diff:
#   1. Reword three existing concept strings (B22, F66, B76).
#   2. Turn the "Partial Derivative Machine Derivatives" activity (row 83)
#      from a placeholder ("[]" representation + auto long-description)
#      into a fully filled-in row: a real [Representations] entry and a
#      hand-written Long Description (replacing the generic formula).
#   3. Nudge the saved view state (freeze-pane scroll position + the
#      active/selected cell) to where the author's cursor ended up.
#
# All the dependent/derived cells (the "[...]" roll-up formulas in columns
# D/E on other rows, and the shared "long description" formula in column I)
# recalculate on their own because they reference these cells by formula.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- 1. Reword existing concept / representation strings --------------

# Row 22, column B ("Name"): sharpen the wording of the concept.
$ws.Cells.Item(22, 2).Value = "The value of a partial derivative depend on the value(s) of what is held constant"

# Row 22 grows from one wrapped line to two once the longer text is in
# place (wrap-text is on for every cell in this sheet); match the height
# Excel's own auto-fit would have produced.
$ws.Rows.Item(22).RowHeight = 47.25

# Row 66, column F ("[Representations]"): pluralize "Contour Map".
$ws.Cells.Item(66, 6).Value = "[Kinesthetic, Vector Field Map, Contour Maps]"

# Row 76, column B ("Name"): clarify the wording.
$ws.Cells.Item(76, 2).Value = "Partial derivatives depend on (are defined by?) what you hold constant"

# --- 2. Fill in the "Partial Derivative Machine Derivatives" row (83) --

# Column F ("[Representations]"): was the empty placeholder "[]".
$ws.Cells.Item(83, 6).Value = "[partial f/partial x rightarrow partial f/partial x fixing y, picture of PDM, data table (pic)]"

# Column I ("Long Description"): replace the generic shared formula with
# a real, hand-written long description for this activity.
$ws.Cells.Item(83, 9).Value = "In this activity, students experimentally determine various derivatives using the partial derivate machine, a mechanical analogue for thermodynamic systems. Students explore the ratio, limit, and function aspects of multi-variable derivatives, with an emphasis on holding different variables constant. This activity is also an excellent exercise in representational fluency, as students must coordinate experiments and tables of data with (new) symbolic notations."

# --- 3. View state: scroll position + active selection -----------------

$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 58
$win.ScrollColumn = 1
$ws.Range("F67").Select()
